$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 'Rv1600'
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = 'hisC hisC1 Rv1600 MTCY336.04c'
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = 19

$ws.Cells.Item(4, 1).Value = 'Rv1408'
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 'rpe Rv1408 MTCY21B4.25'
$ws.Cells.Item(4, 4).Value = 'FUNCTION: Catalyzes the reversible epimerization of D-ribulose 5-phosphate to D-xylulose 5-phosphate. {ECO:0000255|HAMAP-Rule:MF_02227}.'
$ws.Cells.Item(4, 5).Value = 19

$ws.Cells.Item(5, 1).Value = 'Rv2439c'
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = 'proB Rv2439c MTCY428.07'
$ws.Cells.Item(5, 4).Value = 'FUNCTION: Catalyzes the transfer of a phosphate group to glutamate to form L-glutamate 5-phosphate. {ECO:0000255|HAMAP-Rule:MF_00456}.'
$ws.Cells.Item(5, 5).Value = 19

$ws.Cells.Item(6, 1).Value = 'Rv2201'
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = 'asnB Rv2201 MTCY190.12'
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = 19

$ws.Cells.Item(7, 1).Value = 'Rv3805c'
$ws.Cells.Item(7, 2).Value = 3
$ws.Cells.Item(7, 3).Value = 'aftB Rv3805c'
$ws.Cells.Item(7, 4).Value = 'FUNCTION: Involved in the biosynthesis of the arabinogalactan (AG) region of the mycolylarabinogalactan-peptidoglycan (mAGP) complex, an essential component of the mycobacterial cell wall. Catalyzes the transfer of arabinofuranosyl (Araf) residues from the sugar donor decaprenyl-phospho-arabinose (DPA) to the arabinan domain to form terminal beta-(1->2)-linked Araf residues, which marks the end point for AG arabinan biosynthesis before decoration with mycolic acids. {ECO:0000269|PubMed:17387176}.'
$ws.Cells.Item(7, 5).Value = 19

$ws.Cells.Item(8, 1).Value = 'Rv2063'
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = 'mazF7 Rv2063A'
$ws.Cells.Item(8, 4).Value = 'FUNCTION: Toxic component of a type II toxin-antitoxin (TA) system. Upon expression in E.coli and M.smegmatis inhibits cell growth and colony formation. Its toxic effect is neutralized by coexpression with cognate antitoxin MazE7 (PubMed:19016878, PubMed:20011113). Probably an endoribonuclease (By similarity). {ECO:0000250|UniProtKB:P9WIH9, ECO:0000269|PubMed:19016878, ECO:0000269|PubMed:20011113}.'
$ws.Cells.Item(8, 5).Value = 19

$ws.Cells.Item(9, 1).Value = 'Rv0127'
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = 'mak Rv0127'
$ws.Cells.Item(9, 4).Value = 'FUNCTION: Catalyzes the ATP-dependent phosphorylation of maltose to maltose 1-phosphate (By similarity). Is involved in a branched alpha-glucan biosynthetic pathway from trehalose, together with TreS, GlgE and GlgB. {ECO:0000250, ECO:0000269|PubMed:20305657}.'
$ws.Cells.Item(9, 5).Value = 19

$ws.Cells.Item(10, 1).Value = 'Rv0553'
$ws.Cells.Item(10, 2).Value = 3
$ws.Cells.Item(10, 3).Value = 'menC Rv0553 MTCY25D10.32'
$ws.Cells.Item(10, 4).Value = 'FUNCTION: Converts 2-succinyl-6-hydroxy-2,4-cyclohexadiene-1-carboxylate (SHCHC) to 2-succinylbenzoate (OSB). {ECO:0000255|HAMAP-Rule:MF_00470}.'
$ws.Cells.Item(10, 5).Value = 19

$ws.Cells.Item(11, 1).Value = 'Rv2063'
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = 'mazE7 Rv2063'
$ws.Cells.Item(11, 4).Value = 'FUNCTION: Antitoxin component of a type II toxin-antitoxin (TA) system. Upon expression in E.coli but not in M.smegmatis neutralizes the effect of cognate toxin MazF7. {ECO:0000269|PubMed:19016878}.'
$ws.Cells.Item(11, 5).Value = 19

$ws.Cells.Item(12, 1).Value = 'Rv3902c'
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = 'Rv3902c LH57_21250'
$ws.Cells.Item(12, 4).Value = 'FUNCTION: Antitoxin for tuberculosis necrotizing toxin (TNT). Acts by binding directly to TNT, which inhibits NAD(+) glycohydrolase activity of TNT and protects M.tuberculosis from self-poisoning. {ECO:0000269|PubMed:26237511}.'
$ws.Cells.Item(12, 5).Value = 19

$ws.Cells.Item(13, 1).Value = 'Rv1254'
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 'Rv1254'
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = 19

$ws.Cells.Item(14, 1).Value = 'Rv1479'
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 'moxR1 Rv1479'
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = 19

$ws.Cells.Item(15, 1).Value = 'Rv3594'
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(15, 3).Value = 'Rv3594'
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = 19

$ws.Cells.Item(16, 1).Value = 'Rv0378'
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = 'Rv0378'
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = 19
